$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook uses 4 "emoji" marker strings in column A (intervention_type
# status icons). Excel keeps these as shared strings; we need to replace the
# emoji text with the new labels while keeping the cells as plain text
# (not numbers) and not introducing any new cell styles/number formats.
#
# Plain `Range.Value = "-3"` (or "+3") gets auto-coerced by Excel into a
# number, so instead we write a literal-text formula ("="-3"") and then
# convert the whole range to static values in one PasteSpecial (values-only)
# pass; this keeps the result as shared-string text without touching
# NumberFormat/quote-prefix styles.

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

$changed = $false

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $t = $cell.Text
    if ($t -eq "📘") {
        $cell.Formula = "=""⚠️"""
        $changed = $true
    } elseif ($t -eq "📕") {
        $cell.Formula = "=""-3"""
        $changed = $true
    } elseif ($t -eq "📙") {
        $cell.Formula = "=""+3"""
        $changed = $true
    } elseif ($t -eq "📗") {
        $cell.Formula = "=""✅"""
        $changed = $true
    }
}

if ($changed) {
    $rng = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1))
    $rng.Copy()
    $rng.PasteSpecial(-4163)
    $excel.CutCopyMode = 0
}
